# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Re-set header values since PasteSpecial(formats) shouldn't clobber them, but ensure correctness
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in season record data for rows 2-44
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 92   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 71   # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}
